# Update the LR-pairs worksheet with new TPM-derived values and drop the
# now-obsolete MuSCs/Resolving-Mac combination rows (old rows 6-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 through 9 (entire rows) - data no longer present after re-run
$ws.Range("A6:T9").EntireRow.Delete()

# Row 2: ECs -> Cd200 -> Cd200r1 -> Resolving-Mac
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 80.87054666666667
$ws.Range("H2").Value = 242.61164
$ws.Range("I2").Value = 0.7161501349062054
$ws.Range("J2").Value = 0.7161501349062055
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.376254
$ws.Range("N2").Value = 31.128762
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 839.1333333321868
$ws.Range("R2").Value = 7552.199999989681
$ws.Range("S2").Value = 0.7161501349062054
$ws.Range("T2").Value = 0.7161501349062055

# Row 3: FAPs -> Cd200 -> Cd200r1 -> Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 4.561623333333333
$ws.Range("H3").Value = 13.68487
$ws.Range("I3").Value = 0.04039551233681073
$ws.Range("J3").Value = 0.04039551233681073
$ws.Range("M3").Value = 10.376254
$ws.Range("N3").Value = 31.128762
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 47.33256235899334
$ws.Range("R3").Value = 425.99306123094
$ws.Range("S3").Value = 0.04039551233681073
$ws.Range("T3").Value = 0.04039551233681073

# Row 4: MuSCs -> Cd200 -> Cd200r1 -> Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 27.01376833333333
$ws.Range("H4").Value = 81.04130499999999
$ws.Range("I4").Value = 0.2392207624857774
$ws.Range("J4").Value = 0.2392207624857774
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.376254
$ws.Range("N4").Value = 31.128762
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 280.3017217238233
$ws.Range("R4").Value = 2522.71549551441
$ws.Range("S4").Value = 0.2392207624857774
$ws.Range("T4").Value = 0.2392207624857774

# Row 5: Resolving-Mac -> Cd200 -> Cd200r1 -> Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4780740000000001
$ws.Range("H5").Value = 1.434222
$ws.Range("I5").Value = 0.004233590271206475
$ws.Range("J5").Value = 0.004233590271206476
$ws.Range("M5").Value = 10.376254
$ws.Range("N5").Value = 31.128762
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 4.960617254796001
$ws.Range("R5").Value = 44.64555529316401
$ws.Range("S5").Value = 0.004233590271206475
$ws.Range("T5").Value = 0.004233590271206476
